# Two-digit/one-digit division worksheet: refresh the 25 problems in the
# first table (5 rows of 5 cells each, rows 1/5/9/13/17 of the 20-row
# table hold the problems; the rows between are left blank for work).
#
# NOTE: several problems share identical "before" text (e.g. "48÷8="
# appears twice, at row 5 col 2 and row 13 col 3, with different
# replacement targets). Doing a single document-wide Find & Replace
# would only be able to target one value for all matching occurrences,
# and a scoped Find.Execute(..., Replace:=2) was observed to mutate the
# shared underlying text run and bleed into the other occurrence too.
# To keep each substitution local to its own cell we locate the text
# with Find.Execute (no replacement argument) scoped to that cell's
# Range, then overwrite that cell's Range.Text directly.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
$c.Range.Find.Execute("70÷2=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "29÷8="
$c = $t.Cell(1, 2)
$c.Range.Find.Execute("83÷3=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "83÷7="
$c = $t.Cell(1, 3)
$c.Range.Find.Execute("79÷2=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "26÷6="
$c = $t.Cell(1, 4)
$c.Range.Find.Execute("65÷5=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "27÷5="
$c = $t.Cell(1, 5)
$c.Range.Find.Execute("46÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "14÷9="

$c = $t.Cell(5, 1)
$c.Range.Find.Execute("72÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "92÷2="
$c = $t.Cell(5, 2)
$c.Range.Find.Execute("48÷8=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "90÷9="
$c = $t.Cell(5, 3)
$c.Range.Find.Execute("58÷8=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "39÷3="
$c = $t.Cell(5, 4)
$c.Range.Find.Execute("85÷5=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "61÷5="
$c = $t.Cell(5, 5)
$c.Range.Find.Execute("68÷7=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "12÷4="

$c = $t.Cell(9, 1)
$c.Range.Find.Execute("59÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "27÷5="
$c = $t.Cell(9, 2)
$c.Range.Find.Execute("29÷5=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "58÷2="
$c = $t.Cell(9, 3)
$c.Range.Find.Execute("92÷7=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "85÷2="
$c = $t.Cell(9, 4)
$c.Range.Find.Execute("75÷7=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "79÷9="
$c = $t.Cell(9, 5)
$c.Range.Find.Execute("51÷4=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "91÷3="

$c = $t.Cell(13, 1)
$c.Range.Find.Execute("39÷2=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "51÷2="
$c = $t.Cell(13, 2)
$c.Range.Find.Execute("18÷3=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "86÷8="
$c = $t.Cell(13, 3)
$c.Range.Find.Execute("48÷8=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "13÷9="
$c = $t.Cell(13, 4)
$c.Range.Find.Execute("95÷4=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "47÷7="
$c = $t.Cell(13, 5)
$c.Range.Find.Execute("43÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "83÷7="

$c = $t.Cell(17, 1)
$c.Range.Find.Execute("24÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "82÷2="
$c = $t.Cell(17, 2)
$c.Range.Find.Execute("93÷6=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "85÷2="
$c = $t.Cell(17, 3)
$c.Range.Find.Execute("11÷2=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "18÷6="
$c = $t.Cell(17, 4)
$c.Range.Find.Execute("50÷2=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "42÷6="
$c = $t.Cell(17, 5)
$c.Range.Find.Execute("84÷4=", $false, $false, $false, $false, $false, $true, 1, $false) | Out-Null
$c.Range.Text = "65÷9="
